# Delay.pptx edits:
#  1) Refresh the cached "datetimeFigureOut" header/footer date field text
#     (07/12/2023 -> 21/12/2023) everywhere it is cached: the slide master,
#     every slide layout that carries a Date placeholder, and the notes
#     master.
#  2) Slide 4 ("Delay with feedback"): widen/reposition the body
#     placeholder (Rectangle 3) so it starts further left and is wider.

$EMU_PER_POINT = 12700

function Update-DateField {
    param($shapes, $oldText, $newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq $oldText) {
                    $tr.Text = $newText
                }
            }
        }
    }
}

$p = $ppt.ActivePresentation

$oldDate = "07/12/2023"
$newDate = "21/12/2023"

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DateField $master.Shapes $oldDate $newDate

# Every slide layout hanging off the slide master (Date placeholders live
# on whichever layouts carry one; layouts without one are no-ops).
foreach ($layout in $master.CustomLayouts) {
    Update-DateField $layout.Shapes $oldDate $newDate
}

# Notes master date placeholder.
$notesMaster = $p.NotesMaster
Update-DateField $notesMaster.Shapes $oldDate $newDate

# Slide 4: "Delay with feedback" - resize/reposition the body placeholder.
$slide4 = $p.Slides.Item(4)
$bodyShape = $slide4.Shapes.Item(2)
$bodyShape.Left = 398454 / $EMU_PER_POINT
$bodyShape.Width = 10269546 / $EMU_PER_POINT
$bodyShape.Top = 765175 / $EMU_PER_POINT
$bodyShape.Height = 3384550 / $EMU_PER_POINT
